# Add a new "2017_monthly" worksheet at the end of the workbook containing
# monthly resident / non-resident termination counts for 2017, and make it
# the active/selected sheet (matching the authored commit).

$wb = $excel.ActiveWorkbook

# --- Create the new worksheet as the last tab --------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "2017_monthly"

# --- Data ----------------------------------------------------------------------
$months      = @("January","February","March","April","May","June","July","August","September","October","November","December")
$resident    = @(571,617,741,584,588,554,513,632,588,580,615,589)
$nonresident = @(47,49,51,42,55,66,42,51,47,54,40,62)

# Column headers first (Resident / Non-Resident totals with counts embedded
# as a second line, matching the source report layout).
$ws.Range("B1").Value = "Resident Terminations`r(n = 7,172)"
$ws.Range("C1").Value = "Non-Resident Terminations`r(n = 606)"

# Month rows.
for ($i = 0; $i -lt 12; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $months[$i]
    $ws.Cells.Item($row, 2).Value = $resident[$i]
    $ws.Cells.Item($row, 3).Value = $nonresident[$i]
}

# Row-label header last.
$ws.Range("A1").Value = "month"

# --- Make the new sheet the active / selected tab -----------------------------
$ws.Activate()
$ws.Range("E7").Select()
